# Update specific values in column C to new imputed results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -12.71529999999999
    10 = -13.367
    12 = -10.7193
    18 = -12.2512
    37 = -13.0138
    55 = -13.47429999999999
    68 = -11.496
    77 = -12.1084
    78 = -12.24780000000001
    81 = -13.371
    82 = -12.1584
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
